$wb = $excel.ActiveWorkbook

# --- settings sheet: bump form version/title from V3 to V4 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "ci_lf_pretas_2_participant_202105_v4"
$ws3.Range("A2").Value = "(May 2021) 2. Côte d'Ivoire -  Pre TAS FL Formulaire Participants V4"

# --- survey sheet: widen p_num (order number) constraint from 300 to 500 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D10").Value = "Must be between 1 and 500"
$ws1.Range("F10").Value = "Doit être compris entre 1 et 500 "
$ws1.Range("H10").Value = ". > 0 and . <= 500"
$ws1.Range("I10").Value = "The value must be between 1 and 500"
$ws1.Range("J10").Value = "La valeur doit être compris entre 1 et 500 "
$ws1.Range("E10").Value = "Répeter le numéro d'ordre"

# --- restore survey sheet as the active tab/selection ---
$ws3.Range("A2").Select()
$ws1.Activate()
$ws1.Range("E10").Select()
